# Apply template data changes to the DATA worksheet (sheet "DATA"),
# setting a selection of cells to 0, and updating the active selection
# cell to H3 (matching the author's last-saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

$cellsToZero = @(
    "H3", "I3", "P3",
    "N4", "P4",
    "E5", "F5", "H5", "J5",
    "K7", "L7", "O7",
    "H8", "I8", "N8",
    "J9",
    "H11", "O11",
    "F12",
    "E13", "J13",
    "F14", "H14", "M14",
    "E16", "P16",
    "F17",
    "H18",
    "L20"
)

foreach ($cellRef in $cellsToZero) {
    $ws.Range($cellRef).Value = 0
}

# Update the saved selection/active cell to H3, as in the diff.
$ws.Range("H3").Select()
